# "Table List.xlsx" -> room status table, entered twice (rows 1-3 and rows 9-11).
# Block layout:
#   row 1 : merged "ROOMS" banner, centered
#   row 2 : "ROOM NO." | "STATUS" column headers, centered
#   row 3 : 1 | "occupied" (or "available" in the 2nd block), centered

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- values --------------------------------------------------------------
# Entered in this order so the shared-string table is built as
# ROOMS, occupied, STATUS, "ROOM NO.", available.
$ws.Range("A1").Value = "ROOMS"
$ws.Range("B3").Value = "occupied"
$ws.Range("B2").Value = "STATUS"
$ws.Range("A2").Value = "ROOM NO."
$ws.Range("A3").Value = 1

$ws.Range("A9").Value = "ROOMS"
$ws.Range("A10").Value = "ROOM NO."
$ws.Range("B10").Value = "STATUS"
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "available"

# -- formatting ------------------------------------------------------------
# Center-align every cell across the 3-column block in both copies.
$ws.Range("A1:C1").HorizontalAlignment = -4108
$ws.Range("A2:C3").HorizontalAlignment = -4108
$ws.Range("A9:C9").HorizontalAlignment = -4108
$ws.Range("A10:C11").HorizontalAlignment = -4108

# Merge the "ROOMS" banner across the 3 columns in both blocks.
$ws.Range("A1:C1").Merge() | Out-Null
$ws.Range("A9:C9").Merge() | Out-Null

# Column widths sized to fit the header text.
$ws.Columns("A").ColumnWidth = 9.666666666666666
$ws.Columns("B").ColumnWidth = 8.166666666666666

# Leave the selection on K8, matching the saved workbook state.
$ws.Range("K8").Select() | Out-Null

Write-Host "done"
